$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "Texas Notes" between "Calculations" and
#    "EoDSDwSP".
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsTexas = $wb.Worksheets.Add($null, $wsCalc)
$wsTexas.Name = "Texas Notes"

# Re-resolve sheet handles now that the sheet collection has shifted -
# worksheet references captured before the insert no longer track by name.
$wsAbout = $wb.Worksheets.Item("About")
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsTexas = $wb.Worksheets.Item("Texas Notes")
$wsEoDS = $wb.Worksheets.Item("EoDSDwSP")

# ---------------------------------------------------------------------------
# 2. Populate "Texas Notes" content. The layout mirrors Calculations!A1:D10
#    (Type / BAU deployment / Extended ITC / % increase, ITC incentive
#    level, elasticity), just shifted down four rows, so reuse its formats.
#    (String cells are populated in the original authoring order so the
#    shared-string table indices line up with the source workbook.)
# ---------------------------------------------------------------------------
$wsTexas.Range("B5").Value = "BAU Deployment 2015-2022 (square pixels measured)"
$wsTexas.Range("B5").Font.Bold = $true
$wsTexas.Range("B5").HorizontalAlignment = -4152

$wsTexas.Range("C5").Value = "Deployment with Extended ITC 2015-2022 (square pixels measured)"
$wsTexas.Range("C5").Font.Bold = $true
$wsTexas.Range("C5").HorizontalAlignment = -4152

$wsTexas.Range("A1").Value = "The source has Texas specific data, but it's in graphical form. "
$wsTexas.Range("A2").Value = "So, I used some visual editing software to measure the areas of the Texas specific data. This has some error associated with it."

$wsTexas.Range("A5").Value = "Type"
$wsTexas.Range("A5").Font.Bold = $true

$wsTexas.Range("D5").Value = "% Increase due to ITC"
$wsTexas.Range("D5").Font.Bold = $true
$wsTexas.Range("D5").HorizontalAlignment = -4152

$wsTexas.Range("A6").Value = "Residential"
$wsTexas.Range("B6").Value = 26477
$wsTexas.Range("C6").Value = 29137
$wsTexas.Range("D6").Formula = "=(C6-B6)/B6"
$wsTexas.Range("D6").NumberFormat = "0.00%"

$wsTexas.Range("A7").Value = "Commercial"
$wsTexas.Range("B7").Value = 11010
$wsTexas.Range("C7").Value = 20634
$wsTexas.Range("D7").Formula = "=(C7-B7)/B7"
$wsTexas.Range("D7").NumberFormat = "0.00%"

$wsTexas.Range("A9").Value = "ITC Incentive Level"
$wsTexas.Range("A9").Font.Bold = $true
$wsTexas.Range("A9").HorizontalAlignment = -4152

$wsTexas.Range("A10").Value = 0.3
$wsTexas.Range("A10").NumberFormat = "0%"
$wsTexas.Range("B10").Value = "of system cost"

$wsTexas.Range("A12").Value = "Elasticity of Distributed Solar Deployment with respect to ITC Incentive Level"
$wsTexas.Range("A12").Font.Bold = $true

$wsTexas.Range("A13").Value = "Residential"
$wsTexas.Range("B13").Formula = "=D6/`$A`$10"
$wsTexas.Range("B13").NumberFormat = "0.000"

$wsTexas.Range("A14").Value = "Commercial"
$wsTexas.Range("B14").Formula = "=D7/`$A`$10"
$wsTexas.Range("B14").NumberFormat = "0.000"

$wsTexas.Range("A16").Value = "The main point here is that Texas residential is less elastic and commericial is more elastic than the national average. Because the measuring technique I used has some error in it, I will average the numbres"
$wsTexas.Range("A17").Value = "above with the national numbers in the " + [char]34 + "Calculations" + [char]34 + " tab to come up with something a bit more conservative, in case my measuring error is high."

$wsTexas.Range("A19").Value = "Residential"
$wsTexas.Range("B19").Formula = "=AVERAGE(B13,Calculations!B9)"
$wsTexas.Range("B19").NumberFormat = "0.000"

$wsTexas.Range("A20").Value = "Commercial"
$wsTexas.Range("B20").Formula = "=AVERAGE(B14,Calculations!B10)"
$wsTexas.Range("B20").NumberFormat = "0.000"

$wsTexas.Range("B21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-point EoDSDwSP's elasticity formulas at the new "Texas Notes" sheet.
# ---------------------------------------------------------------------------
$wsEoDS.Range("B2").Formula = "='Texas Notes'!B19"
$wsEoDS.Range("B4").Formula = "='Texas Notes'!B20"
$wsEoDS.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. About sheet: add hyperlink on B6, tidy up selection.
# ---------------------------------------------------------------------------
$wsAbout.Hyperlinks.Add($wsAbout.Range("B6"), "http://www.seia.org/sites/default/files/resources/BNEF_SEIA%20Solar%20Forecast_15%20September%202015.pdf") | Out-Null
$wsAbout.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Calculations sheet selection.
# ---------------------------------------------------------------------------
$wsCalc.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. Make EoDSDwSP the active sheet/tab again (tabSelected).
# ---------------------------------------------------------------------------
$wsEoDS.Activate()
